$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: 2020-04-02, 11:15 -> 24:00 (midnight)
$ws.Cells.Item(12, 2).Value2 = 6977
$ws.Cells.Item(12, 3).Value2 = 43923
$ws.Cells.Item(12, 4).Value2 = 0.46875
$ws.Cells.Item(12, 5).Value2 = 0
$ws.Cells.Item(12, 7).Value2 = "Created a draft of LogicUnit.vhd (have not compiled it yet)"

# Row 13: 2020-04-03, 24:00 -> 00:30
$ws.Cells.Item(13, 2).Value2 = 6977
$ws.Cells.Item(13, 3).Value2 = 43924
$ws.Cells.Item(13, 4).Value2 = 1
$ws.Cells.Item(13, 5).Value2 = 0.020833333333333332
$ws.Cells.Item(13, 7).Value2 = "Copied over Lab 2 into Adder.vhd (have not compiled it yet)"

# Row 14: 2020-04-03, 00:30 -> 00:45
$ws.Cells.Item(14, 2).Value2 = 6977
$ws.Cells.Item(14, 3).Value2 = 43924
$ws.Cells.Item(14, 4).Value2 = 0.020833333333333332
$ws.Cells.Item(14, 5).Value2 = 0.03125
$ws.Cells.Item(14, 7).Value2 = "Formatted LogicUnit.vhd and Adder.vhd to look nicer"

# Update the active selection to A16, matching the saved view state
$ws.Range("A16").Select()
